# Fixed wavelengths not appearing on absorbance worksheets & data frames.
#
# - "Materials" / "Samples" sheets: drop the trailing wavelength row (row 7,
#   wavelength 1000.111) which had no real data.
# - "Standard Deviations": the Wavelength column held real wavelengths but
#   should instead be the plain row index (0..4); also drop trailing row 7.
# - "Absorbance Sample" / "Absorbance Material": these were missing the
#   Wavelength column entirely - insert it as column A (shifting all other
#   columns right by one) and populate it with the wavelengths; also drop
#   the trailing row 7.

$wb = $excel.ActiveWorkbook

$wavelengths = @(630.188, 710.104, 800.131, 905.029, 940.061)

# ---- Materials (sheet1): just remove row 7 ----
$wsMaterials = $wb.Worksheets.Item("Materials")
$wsMaterials.Rows.Item(7).Delete()

# ---- Samples (sheet2): just remove row 7 ----
$wsSamples = $wb.Worksheets.Item("Samples")
$wsSamples.Rows.Item(7).Delete()

# ---- Standard Deviations (sheet3): column A becomes row index, drop row 7 ----
$wsStd = $wb.Worksheets.Item("Standard Deviations")
for ($i = 0; $i -lt 5; $i++) {
    $wsStd.Cells.Item($i + 2, 1).Value = $i
}
$wsStd.Rows.Item(7).Delete()

# ---- Absorbance Sample (sheet4): insert Wavelength column, drop row 7 ----
$wsAbsSample = $wb.Worksheets.Item("Absorbance Sample")
$wsAbsSample.Columns.Item(1).Insert()
$wsAbsSample.Cells.Item(1, 1).Value = "Wavelength"
$wsAbsSample.Cells.Item(1, 2).Copy()
$wsAbsSample.Cells.Item(1, 1).PasteSpecial(-4122)
for ($i = 0; $i -lt 5; $i++) {
    $wsAbsSample.Cells.Item($i + 2, 1).Value = $wavelengths[$i]
}
$wsAbsSample.Rows.Item(7).Delete()

# ---- Absorbance Material (sheet5): insert Wavelength column, drop row 7 ----
$wsAbsMaterial = $wb.Worksheets.Item("Absorbance Material")
$wsAbsMaterial.Columns.Item(1).Insert()
$wsAbsMaterial.Cells.Item(1, 1).Value = "Wavelength"
$wsAbsMaterial.Cells.Item(1, 2).Copy()
$wsAbsMaterial.Cells.Item(1, 1).PasteSpecial(-4122)
for ($i = 0; $i -lt 5; $i++) {
    $wsAbsMaterial.Cells.Item($i + 2, 1).Value = $wavelengths[$i]
}
$wsAbsMaterial.Rows.Item(7).Delete()
